$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "37.206.69"

$ws.Range("D3").Value = "2.067.90"
$ws.Range("E3").Value = "  +1.25%  "

$ws.Range("E4").Value = "  +0.03%  "

$ws.Range("D5").Value = "'249.32"
$ws.Range("E5").Value = "  +0.71%  "

$ws.Range("D6").Value = "'0.668"
$ws.Range("E6").Value = "  +1.14%  "

$ws.Range("D7").Value = "'58.61"
$ws.Range("E7").Value = "  +5.33%  "

$ws.Range("E8").Value = "  -0.09%  "

$ws.Range("E9").Value = "  +2.58%  "

$ws.Range("D10").Value = "'0.0790"
$ws.Range("E10").Value = "  +1.79%  "

$ws.Range("E11").Value = "  +2.20%  "

$ws.Range("D12").Value = "'15.89"
$ws.Range("E12").Value = "  +1.46%  "

$ws.Range("D13").Value = "'0.919"
$ws.Range("E13").Value = "  +17.07%  "

$ws.Range("D14").Value = "2.369.44"
$ws.Range("E14").Value = "  +1.22%  "

$ws.Range("D15").Value = "'5.88"
$ws.Range("E15").Value = "  +4.69%  "

$ws.Range("D16").Value = "2.071.14"
$ws.Range("E16").Value = "  +1.39%  "

$ws.Range("D17").Value = "'18.60"
$ws.Range("E17").Value = "  +13.91%  "

$ws.Range("D18").Value = "37.232.12"
$ws.Range("E18").Value = "  +0.86%  "

$ws.Range("D19").Value = "'75.52"
$ws.Range("E19").Value = "  +2.81%  "

$ws.Range("D20").Value = "0.0₃0912"
$ws.Range("E20").Value = "  +2.51%  "

$ws.Range("D21").Value = "'5.54"
$ws.Range("E21").Value = "  +4.99%  "

$ws.Range("D22").Value = "'239.24"
$ws.Range("E22").Value = "  +1.79%  "

$ws.Range("D23").Value = "'1.00"
$ws.Range("E23").Value = "  -0.10%  "

$ws.Range("E24").Value = "  +6.02%  "

$ws.Range("D25").Value = "'2.23"
$ws.Range("E25").Value = "  +2.59%  "

$ws.Range("D26").Value = "'9.68"
$ws.Range("E26").Value = "  +7.31%  "

$ws.Range("D27").Value = "'171.74"
$ws.Range("E27").Value = "  +2.76%  "

$ws.Range("D28").Value = "'20.33"
$ws.Range("E28").Value = "  +3.74%  "

$ws.Range("D29").Value = "'5.62"
$ws.Range("E29").Value = "  +21.51%  "

$ws.Range("E30").Value = "  +1.56%  "

$ws.Range("D31").Value = "'1.16"
$ws.Range("E31").Value = "  +7.06%  "

$ws.Range("D32").Value = "'4.90"
$ws.Range("E32").Value = "  +11.88%  "

$ws.Range("E33").Value = "  +4.07%  "

$ws.Range("B34").Value = "LidoDAOToken"
$ws.Range("C34").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D34").Value = "'2.36"
$ws.Range("E34").Value = "  +7.59%  "

$ws.Range("B35").Value = "Kaspa"
$ws.Range("C35").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D35").Value = "'0.0884"
$ws.Range("E35").Value = "  +1.60%  "

$ws.Range("E36").Value = "  +0.01%  "

$ws.Range("E37").Value = "  +4.92%  "

$ws.Range("E38").Value = "  +1.66%  "

$ws.Range("D39").Value = "'5.21"
$ws.Range("E39").Value = "  +7.21%  "

$ws.Range("E40").Value = "  -2.87%  "

$ws.Range("E41").Value = "  -3.55%  "

$ws.Range("E42").Value = "  +3.76%  "

$ws.Range("B43").Value = "ARBITRUM"
$ws.Range("C43").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D43").Value = "'1.18"
$ws.Range("E43").Value = "  +6.43%  "

$ws.Range("B44").Value = "Aave"
$ws.Range("C44").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D44").Value = "'100.74"
$ws.Range("E44").Value = "  +6.34%  "

$ws.Range("B45").Value = "InjectiveProtocol"
$ws.Range("C45").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D45").Value = "'17.59"
$ws.Range("E45").Value = "  +3.18%  "

$ws.Range("E46").Value = "  +1.75%  "

$ws.Range("D47").Value = "'3.94"
$ws.Range("E47").Value = "  +19.12%  "

$ws.Range("D48").Value = "1.316.29"
$ws.Range("E48").Value = "  +3.73%  "

$ws.Range("E49").Value = "  +5.66%  "

$ws.Range("D50").Value = "'2.89"
$ws.Range("E50").Value = "  +1.82%  "

$ws.Range("D51").Value = "2.257.61"
$ws.Range("E51").Value = "  +1.41%  "
